# Auto-generated Excel COM-interop script
# Applies cached-value corrections to the Leve profit columns
# (currentAveragePrice / NQ / HQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 7
$ws.Range("H7").Value = 5500
$ws.Range("J7").Value = 5500
$ws.Range("L7").Value = 5500
$ws.Range("N7").Value = -5724
# Row 14
$ws.Range("H14").Value = 5500
$ws.Range("J14").Value = 5500
$ws.Range("L14").Value = 5500
$ws.Range("N14").Value = -5882
# Row 28
$ws.Range("H28").Value = 2612.3
$ws.Range("I28").Value = 1777.9333
$ws.Range("K28").Value = 1777.9333
$ws.Range("M28").Value = -1292.9333
# Row 100
$ws.Range("H100").Value = 1184
$ws.Range("I100").Value = 593
$ws.Range("J100").Value = 1775
$ws.Range("K100").Value = 593
$ws.Range("L100").Value = 1775
$ws.Range("M100").Value = -52
$ws.Range("N100").Value = -2857
# Row 106
$ws.Range("H106").Value = 6690.1875
$ws.Range("I106").Value = 6690.1875
$ws.Range("K106").Value = 6690.1875
$ws.Range("M106").Value = -6059.1875

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4934.3447
$ws.Range("I61").Value = 3056.5217
$ws.Range("K61").Value = 3056.5217
$ws.Range("M61").Value = -2844.5217
# Row 74
$ws.Range("H74").Value = 2111.9333
$ws.Range("I74").Value = 2162.8215
$ws.Range("J74").Value = 1399.5
$ws.Range("K74").Value = 2162.8215
$ws.Range("L74").Value = 1399.5
$ws.Range("M74").Value = -1288.8215
$ws.Range("N74").Value = -3147.5
# Row 77
$ws.Range("H77").Value = 2111.9333
$ws.Range("I77").Value = 2162.8215
$ws.Range("J77").Value = 1399.5
$ws.Range("K77").Value = 10814.1075
$ws.Range("L77").Value = 6997.5
$ws.Range("M77").Value = -6446.1075
$ws.Range("N77").Value = -15733.5
# Row 122
$ws.Range("H122").Value = 5112.7085
$ws.Range("I122").Value = 3982.4119
$ws.Range("K122").Value = 11947.2357
$ws.Range("M122").Value = -9497.235700000001
# Row 136
$ws.Range("H136").Value = 4934.3447
$ws.Range("I136").Value = 3056.5217
$ws.Range("K136").Value = 9169.5651
$ws.Range("M136").Value = -6619.5651

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 25247.762
$ws.Range("I94").Value = 1284.6428
$ws.Range("J94").Value = 73174
$ws.Range("K94").Value = 1284.6428
$ws.Range("L94").Value = 73174
$ws.Range("M94").Value = -833.6428000000001
$ws.Range("N94").Value = -74076
# Row 107
$ws.Range("H107").Value = 1195.9445
$ws.Range("I107").Value = 1077.7142
$ws.Range("K107").Value = 1077.7142
$ws.Range("M107").Value = 842.2858000000001
# Row 134
$ws.Range("H134").Value = 4865.4375
$ws.Range("I134").Value = 4832
$ws.Range("K134").Value = 14496
$ws.Range("M134").Value = -11961

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 50174.832
$ws.Range("I16").Value = 32210
$ws.Range("K16").Value = 32210
$ws.Range("M16").Value = -31923
# Row 58
$ws.Range("H58").Value = 1910.3684
$ws.Range("I58").Value = 1689.5385
$ws.Range("K58").Value = 1689.5385
$ws.Range("M58").Value = -1486.5385
# Row 105
$ws.Range("H105").Value = 3802.0715
$ws.Range("I105").Value = 2702.5
$ws.Range("J105").Value = 10399.5
$ws.Range("K105").Value = 2702.5
$ws.Range("L105").Value = 10399.5
$ws.Range("M105").Value = -955.5
$ws.Range("N105").Value = -13893.5
# Row 113
$ws.Range("H113").Value = 50174.832
$ws.Range("I113").Value = 32210
$ws.Range("K113").Value = 32210
$ws.Range("M113").Value = -30040
# Row 132
$ws.Range("H132").Value = 1446.3793
$ws.Range("I132").Value = 1378.5
$ws.Range("J132").Value = 1659.7142
$ws.Range("K132").Value = 4135.5
$ws.Range("L132").Value = 4979.142599999999
$ws.Range("M132").Value = -1605.5
$ws.Range("N132").Value = -10039.1426
# Row 134
$ws.Range("H134").Value = 1155.125
$ws.Range("I134").Value = 934
$ws.Range("J134").Value = 1523.6666
$ws.Range("K134").Value = 2802
$ws.Range("L134").Value = 4570.9998
$ws.Range("M134").Value = -267
$ws.Range("N134").Value = -9640.9998
# Row 136
$ws.Range("H136").Value = 1910.3684
$ws.Range("I136").Value = 1689.5385
$ws.Range("K136").Value = 5068.6155
$ws.Range("M136").Value = -2518.6155

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 55
$ws.Range("H55").Value = 4840.294
$ws.Range("J55").Value = 7649.1
$ws.Range("L55").Value = 22947.3
$ws.Range("N55").Value = -23301.3
# Row 69
$ws.Range("H69").Value = 1966.3334
$ws.Range("I69").Value = 1966.3334
$ws.Range("K69").Value = 5899.0002
$ws.Range("M69").Value = -5088.0002
# Row 72
$ws.Range("H72").Value = 1966.3334
$ws.Range("I72").Value = 1966.3334
$ws.Range("K72").Value = 17697.0006
$ws.Range("M72").Value = -13641.0006
# Row 74
$ws.Range("H74").Value = 10000
$ws.Range("J74").Value = 10000
$ws.Range("L74").Value = 30000
$ws.Range("N74").Value = -32122
# Row 77
$ws.Range("H77").Value = 10000
$ws.Range("J77").Value = 10000
$ws.Range("L77").Value = 90000
$ws.Range("N77").Value = -100608
# Row 88
$ws.Range("H88").Value = 5802.25
$ws.Range("I88").Value = 7006
$ws.Range("J88").Value = 4598.5
$ws.Range("K88").Value = 21018
$ws.Range("L88").Value = 13795.5
$ws.Range("M88").Value = -20590
$ws.Range("N88").Value = -14651.5
# Row 91
$ws.Range("H91").Value = 5802.25
$ws.Range("I91").Value = 7006
$ws.Range("J91").Value = 4598.5
$ws.Range("K91").Value = 21018
$ws.Range("L91").Value = 13795.5
$ws.Range("M91").Value = -19536
$ws.Range("N91").Value = -16759.5
# Row 99
$ws.Range("H99").Value = 2097.2
$ws.Range("I99").Value = 2566.3333
$ws.Range("J99").Value = 1393.5
$ws.Range("K99").Value = 7698.999899999999
$ws.Range("L99").Value = 4180.5
$ws.Range("M99").Value = -5452.999899999999
$ws.Range("N99").Value = -8672.5
# Row 140
$ws.Range("H140").Value = 2508.8572
$ws.Range("I140").Value = 1277.1666
$ws.Range("K140").Value = 3831.4998
$ws.Range("M140").Value = 1348.5002
# Row 141
$ws.Range("H141").Value = 4999.375
$ws.Range("I141").Value = 5599
$ws.Range("K141").Value = 16797
$ws.Range("M141").Value = -11617

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 381.2857
$ws.Range("I2").Value = 334
$ws.Range("J2").Value = 499.5
$ws.Range("K2").Value = 334
$ws.Range("L2").Value = 499.5
$ws.Range("M2").Value = -221
$ws.Range("N2").Value = -725.5
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
# Row 97
$ws.Range("H97").Value = 680.6667
$ws.Range("I97").Value = 657.1053000000001
$ws.Range("K97").Value = 657.1053000000001
$ws.Range("M97").Value = -161.1053000000001
# Row 132
$ws.Range("H132").Value = 4870.84
$ws.Range("I132").Value = 5000.9546
$ws.Range("K132").Value = 15002.8638
$ws.Range("M132").Value = -12472.8638

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 799
$ws.Range("I61").Value = 823.875
$ws.Range("K61").Value = 823.875
$ws.Range("M61").Value = -621.875
# Row 113
$ws.Range("H113").Value = 799
$ws.Range("I113").Value = 823.875
$ws.Range("K113").Value = 823.875
$ws.Range("M113").Value = 1346.125

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 63
$ws.Range("H63").Value = 6000
$ws.Range("J63").Value = 6000
$ws.Range("L63").Value = 6000
$ws.Range("N63").Value = -7248
# Row 66
$ws.Range("H66").Value = 6000
$ws.Range("J66").Value = 6000
$ws.Range("L66").Value = 18000
$ws.Range("N66").Value = -24240
# Row 107
$ws.Range("H107").Value = 298
$ws.Range("I107").Value = 298
$ws.Range("K107").Value = 894
$ws.Range("M107").Value = 1026
# Row 126
$ws.Range("H126").Value = 2717.9443
$ws.Range("I126").Value = 2717.9443
$ws.Range("K126").Value = 8153.8329
$ws.Range("M126").Value = -5683.8329
# Row 136
$ws.Range("H136").Value = 3937.9697
$ws.Range("I136").Value = 4094.84
$ws.Range("K136").Value = 12284.52
$ws.Range("M136").Value = -9734.52

Write-Output "Applied Omega_Profits cached-value updates."
